$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "How to advertise for a Voluntary RA" paragraph: replace the reference to
#    a "template attached" with a hyperlink to the MS Forms based
#    "Research Assistant Scheme Project Description Form".
# ---------------------------------------------------------------------------

$oldIntro = "please fill out the template attached and provide a short description of the research project(s) you would like your potential RAs to be working on by"
$newIntro = "please fill out the"
$d.Content.Find.Execute($oldIntro, $true, $false, $false, $false, $false, $true, 1, $false, $newIntro, 2) | Out-Null

# Insert the plain-text placeholder (together with the trailing " by") right
# after "...please fill out the" so that everything is plain text first -
# this keeps the later " by" text from inheriting the hyperlink's character
# style once the placeholder is converted into a real hyperlink.
$insertRange = $d.Content
$insertRange.Find.Execute("please fill out the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertRange.Collapse(0)
$insertRange.InsertAfter(" ~~RASCHEMELINK~~ by")

# Turn the placeholder text into a hyperlink pointing at the online form.
$linkRange = $d.Content
$linkRange.Find.Execute("~~RASCHEMELINK~~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($linkRange, "https://forms.office.com/e/RAProjectDescription", "", "", "Research Assistant Scheme Project Description Form") | Out-Null

# ---------------------------------------------------------------------------
# 2. "Deadlines" list: trim the first bullet so it no longer refers to the
#    mailto link that is being removed, and delete that now-unused bullet
#    paragraph (and its hyperlink) entirely.
# ---------------------------------------------------------------------------

$oldDeadline = ": Project summaries to be sent to Matt using the link below. If you" + [char]0x2019 + "re reading this in Word then control-click on the link."
$newDeadline = ": Project summaries to be sent to Matt"
$d.Content.Find.Execute($oldDeadline, $true, $false, $false, $false, $false, $true, 1, $false, $newDeadline, 2) | Out-Null

# Remove the paragraph that only contained the mailto hyperlink
# (mgreen@bournemouth.ac.uk?subject=Voluntary RA Scheme Application).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "mgreen@bournemouth.ac.uk*") {
        $para.Range.Delete()
        break
    }
}
